# Updated cryptos list with GitHub Actions
# Values that look numeric are prefixed with a leading apostrophe so Excel
# keeps storing them as literal text (matching the sheet's existing
# text-as-string layout) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.544.55"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "'1.722.10"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("D4").Value = "'0.9952"

$ws.Range("D5").Value = "'240.42"
$ws.Range("E5").Value = "  -2.52%  "

$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").Value = "'0.4914"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'0.2590"
$ws.Range("E8").Value = "  -3.41%  "

$ws.Range("D9").Value = "'0.06195"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").Value = "'1.722.50"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'15.69"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.06961"
$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").Value = "'0.6054"
$ws.Range("E13").Value = "  -1.83%  "

$ws.Range("D14").Value = "'4.466"
$ws.Range("E14").Value = "  -2.72%  "

$ws.Range("D15").Value = "'76.64"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "'0.9959"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").Value = "'26.371.26"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "'0.9948"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "'0.000007135"
$ws.Range("E19").Value = "  -2.36%  "

$ws.Range("D20").Value = "'11.33"
$ws.Range("E20").Value = "  -2.25%  "

$ws.Range("D21").Value = "'1.945.14"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").Value = "'4.409"
$ws.Range("E22").Value = "  -3.56%  "

$ws.Range("D23").Value = "'8.436"
$ws.Range("E23").Value = "  -3.41%  "

$ws.Range("D24").Value = "'5.092"
$ws.Range("E24").Value = "  -3.44%  "

$ws.Range("D25").Value = "'137.96"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "'15.24"
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("D28").Value = "'1.747"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").Value = "'105.82"
$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").Value = "'3.912"
$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("D31").Value = "'0.07932"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("D32").Value = "'3.628"
$ws.Range("E32").Value = "  -3.03%  "

$ws.Range("D33").Value = "'0.04485"

$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").Value = "'0.9991"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").Value = "'0.6183"
$ws.Range("E37").Value = "  -3.55%  "

$ws.Range("D38").Value = "'0.9448"
$ws.Range("E38").Value = "  +4.91%  "

$ws.Range("D39").Value = "'2.009"
$ws.Range("E39").Value = "  -3.42%  "

$ws.Range("D40").Value = "'2.389"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("D41").Value = "'0.9957"
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("D42").Value = "'0.01493"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("D44").Value = "'5.473"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "'0.3828"
$ws.Range("E45").Value = "  -2.45%  "

$ws.Range("D46").Value = "'6.941"
$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("D49").Value = "'30.54"
$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").Value = "'7.751"
$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").Value = "'51.39"
$ws.Range("E51").Value = "  -1.02%  "
